$wb = $excel.ActiveWorkbook

# Citywide Totals
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("E2").Value = 75
$ws.Range("B3").Value = 81
$ws.Range("E3").Value = 150
$ws.Range("G3").Value = 148
$ws.Range("J3").Value = 240
$ws.Range("I6").Value = 24
$ws.Range("B9").Value = 388
$ws.Range("C9").Value = 500
$ws.Range("D9").Value = 434
$ws.Range("E9").Value = 501
$ws.Range("F9").Value = 571
$ws.Range("G9").Value = 445
$ws.Range("H9").Value = 472
$ws.Range("K9").Value = 532
$ws.Range("B10").Value = 1401
$ws.Range("C10").Value = 1656
$ws.Range("D10").Value = 1877
$ws.Range("E10").Value = 2295
$ws.Range("F10").Value = 2185
$ws.Range("H10").Value = 630
$ws.Range("K10").Value = 705
$ws.Range("B11").Value = 1932
$ws.Range("C11").Value = 2319
$ws.Range("D11").Value = 2557
$ws.Range("E11").Value = 3035
$ws.Range("F11").Value = 3010
$ws.Range("G11").Value = 1608
$ws.Range("H11").Value = 1394
$ws.Range("I11").Value = 1742
$ws.Range("J11").Value = 1589
$ws.Range("K11").Value = 1646

# By Neighborhood
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("C5").Value = 22
$ws.Range("D5").Value = 20
$ws.Range("F5").Value = 44
$ws.Range("H6").Value = 5
$ws.Range("B7").Value = 56
$ws.Range("D8").Value = 64
$ws.Range("E8").Value = 127
$ws.Range("F10").Value = 21
$ws.Range("H14").Value = 4
$ws.Range("C16").Value = 5
$ws.Range("B19").Value = 51
$ws.Range("D19").Value = 62
$ws.Range("C23").Value = 20
$ws.Range("B28").Value = 106
$ws.Range("C28").Value = 141
$ws.Range("D28").Value = 105
$ws.Range("E28").Value = 94
$ws.Range("F28").Value = 131
$ws.Range("J28").Value = 68
$ws.Range("D29").Value = 28
$ws.Range("K29").Value = 32
$ws.Range("E32").Value = 172
$ws.Range("F32").Value = 200
$ws.Range("G32").Value = 100
$ws.Range("H32").Value = 79
$ws.Range("F36").Value = 91
$ws.Range("D37").Value = 4
$ws.Range("G40").Value = 2
$ws.Range("C44").Value = 2
$ws.Range("E45").Value = 24
$ws.Range("H47").Value = 50
$ws.Range("E48").Value = 15
$ws.Range("G49").Value = 10
$ws.Range("E52").Value = 37
$ws.Range("B53").Value = 276
$ws.Range("C53").Value = 407
$ws.Range("D53").Value = 610
$ws.Range("E53").Value = 778
$ws.Range("G53").Value = 252
$ws.Range("H53").Value = 229
$ws.Range("B61").Value = 21
$ws.Range("D61").Value = 30
$ws.Range("K61").Value = 6
$ws.Range("I68").Value = 12
$ws.Range("F77").Value = 76
$ws.Range("E78").Value = 52
$ws.Range("D80").Value = 32
$ws.Range("H80").Value = 17
$ws.Range("K87").Value = 34
$ws.Range("F88").Value = 9
$ws.Range("D90").Value = 4
$ws.Range("D95").Value = 57
$ws.Range("F95").Value = 66
$ws.Range("H95").Value = 19
$ws.Range("D96").Value = 32
$ws.Range("E96").Value = 37
$ws.Range("F97").Value = 22
$ws.Range("B99").Value = 1932
$ws.Range("C99").Value = 2319
$ws.Range("D99").Value = 2557
$ws.Range("E99").Value = 3035
$ws.Range("F99").Value = 3010
$ws.Range("G99").Value = 1608
$ws.Range("H99").Value = 1394
$ws.Range("I99").Value = 1742
$ws.Range("J99").Value = 1589
$ws.Range("K99").Value = 1646

# Roseland
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("F9").Value = 53
$ws.Range("F10").Value = 76

# Auburn Gresham
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("B5").Value = 19
$ws.Range("B7").Value = 56

# Austin
$ws = $wb.Worksheets.Item('Austin')
$ws.Range("E2").Value = 5
$ws.Range("D8").Value = 32
$ws.Range("D9").Value = 64
$ws.Range("E9").Value = 127

# Washington Heights
$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("F6").Value = 5
$ws.Range("F7").Value = 9

# Garfield Park
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("E7").Value = 55
$ws.Range("F7").Value = 53
$ws.Range("G7").Value = 37
$ws.Range("H7").Value = 39
$ws.Range("E9").Value = 172
$ws.Range("F9").Value = 200
$ws.Range("G9").Value = 100
$ws.Range("H9").Value = 79

# Chatham
$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("D8").Value = 17
$ws.Range("B9").Value = 43
$ws.Range("B10").Value = 51
$ws.Range("D10").Value = 62

# Grand Crossing
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("F9").Value = 38
$ws.Range("F10").Value = 91

# Loop
$ws = $wb.Worksheets.Item('Loop')
$ws.Range("E3").Value = 18
$ws.Range("G3").Value = 25
$ws.Range("C8").Value = 43
$ws.Range("E8").Value = 66
$ws.Range("H8").Value = 78
$ws.Range("B9").Value = 227
$ws.Range("C9").Value = 346
$ws.Range("D9").Value = 534
$ws.Range("E9").Value = 688
$ws.Range("B10").Value = 276
$ws.Range("C10").Value = 407
$ws.Range("D10").Value = 610
$ws.Range("E10").Value = 778
$ws.Range("G10").Value = 252
$ws.Range("H10").Value = 229

# Armour Square
$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("C5").Value = 6
$ws.Range("F5").Value = 12
$ws.Range("D6").Value = 13
$ws.Range("C7").Value = 22
$ws.Range("D7").Value = 20
$ws.Range("F7").Value = 44

# Sheffield & DePaul
$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range("D5").Value = 4
$ws.Range("H5").Value = 3
$ws.Range("D7").Value = 32
$ws.Range("H7").Value = 17

# Uptown
$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K7").Value = 11
$ws.Range("K9").Value = 34

# Rush & Division
$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("E4").Value = 6
$ws.Range("E6").Value = 52

# Englewood
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("B3").Value = 5
$ws.Range("J3").Value = 11
$ws.Range("B8").Value = 68
$ws.Range("C8").Value = 91
$ws.Range("D8").Value = 57
$ws.Range("E8").Value = 63
$ws.Range("F8").Value = 71
$ws.Range("B9").Value = 106
$ws.Range("C9").Value = 141
$ws.Range("D9").Value = 105
$ws.Range("E9").Value = 94
$ws.Range("F9").Value = 131
$ws.Range("J9").Value = 68

# Lake View
$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("H6").Value = 16
$ws.Range("H8").Value = 50

# Jefferson Park
$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("E6").Value = 22
$ws.Range("E7").Value = 24

# Fuller Park
$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("K7").Value = 20
$ws.Range("D8").Value = 21
$ws.Range("D9").Value = 28
$ws.Range("K9").Value = 32

# Woodlawn
$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("F6").Value = 15
$ws.Range("F7").Value = 22

# Logan Square
$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("E7").Value = 31
$ws.Range("E8").Value = 37

# Douglas
$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("C6").Value = 5
$ws.Range("C8").Value = 20

# Bridgeport
$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("H5").Value = 1
$ws.Range("H7").Value = 4

# Lincoln Square
$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("G3").Value = 1
$ws.Range("G7").Value = 10

# West Town
$ws = $wb.Worksheets.Item('West Town')
$ws.Range("H5").Value = 6
$ws.Range("D6").Value = 45
$ws.Range("F6").Value = 58
$ws.Range("D7").Value = 57
$ws.Range("F7").Value = 66
$ws.Range("H7").Value = 19

# Wicker Park
$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("D7").Value = 29
$ws.Range("E7").Value = 31
$ws.Range("D8").Value = 32
$ws.Range("E8").Value = 37

# Lincoln Park
$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("E6").Value = 9
$ws.Range("E7").Value = 15

# Jackson Park
$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 2

# O'Hare
$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("I4").Value = 1
$ws.Range("I7").Value = 12

# Avondale
$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("F6").Value = 16
$ws.Range("F7").Value = 21

# Bucktown
$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("C5").Value = 4
$ws.Range("C6").Value = 5

# Grant Park
$ws = $wb.Worksheets.Item('Grant Park')
$ws.Range("D5").Value = 4
$ws.Range("D6").Value = 4

# West Elsdon
$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("D4").Value = 2
$ws.Range("D5").Value = 4

# Hermosa
$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("E2").Value = 1
$ws.Range("E5").Value = 2

# Ashburn
$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("H5").Value = 3
$ws.Range("H6").Value = 5
